$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I5").Value = -0.6830381338146789
$ws.Range("J5").Value = 0.4528597409125978
$ws.Range("K5").Value = 0.2594084689748836
$ws.Range("L5").Value = 2.655278111889562
